$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Date: 2025-12-26T15:22:58+00:00 -> 2026-01-01T13:37:23+00:00
$meta.Range("B8").Value = "2026-01-01T13:37:23+00:00"

# Description: drop the trailing "Supports goal-directed care planning and
# intervention tracking." sentence. This shared string is also used for the
# "Definition" cell of the first element row (Elements!M2), so updating the
# text here keeps both cells in sync with the new wording.
$newDescription = "Extension to link nursing interventions to the patient goals they are intended to achieve."
$meta.Range("B11").Value = $newDescription
$elements.Range("M2").Value = $newDescription

# Extension.value[x] Type(s): onc-patient-goal -> onc-nursing-goal (keep the
# trailing newline / preserved whitespace from the source cell).
$elements.Range("K6").Value = "Reference(https://clinyqai.github.io/open-nursing-core-ig/StructureDefinition/onc-nursing-goal)
"

# Column K width: 78.00390625 -> 78.18359375 (character units). The COM
# ColumnWidth setter snaps to the nearest on-screen pixel (1/6 of a
# character at this font/size), so 77.3 lands on the closest reachable
# stored width to the target.
$elements.Columns.Item(11).ColumnWidth = 77.3
